# Updated cryptos list on Sat Oct 19 04:43:38 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns on Sheet1 with the
# latest scraped values. Numeric-looking Price values are entered with a
# leading apostrophe so Excel keeps them as text (matching the source
# data, e.g. "598.41" / "0.0000190") instead of auto-converting them to
# numbers; the quote-prefix formatting that introduces is then reset back
# to the default "Normal" style so no stray number format is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.402.35'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = '2.644.80'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''598.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').Value = '''154.43'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.18%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D9').Value = '2.644.90'
$ws.Range('E9').Value = '  +1.06%  '
$ws.Range('E10').Value = '  +8.07%  '
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('D12').Value = '''5.26'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('D13').Value = '''0.354'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.21%  '
$ws.Range('D14').Value = '''28.18'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.24%  '
$ws.Range('D15').Value = '''0.0000190'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.23%  '
$ws.Range('D16').Value = '3.123.39'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').Value = '68.347.26'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').Value = '2.658.76'
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('D19').Value = '''11.45'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.16%  '
$ws.Range('D20').Value = '''365.94'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.67%  '
$ws.Range('E21').Value = '  +13.08%  '
$ws.Range('E22').Value = '  +3.65%  '
$ws.Range('D23').Value = '''4.89'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.81%  '
$ws.Range('D24').Value = '''2.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('D25').Value = '''73.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '''9.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('E28').Value = '  +2.58%  '
$ws.Range('D29').Value = '2.779.46'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').Value = '''575.68'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.92%  '
$ws.Range('D32').Value = '''8.18'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.84%  '
$ws.Range('E33').Value = '  +3.30%  '
$ws.Range('E34').Value = '  +1.40%  '
$ws.Range('D35').Value = '''0.131'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.25%  '
$ws.Range('E36').Value = '  +5.67%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').Value = '''160.35'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.31%  '
$ws.Range('D39').Value = '''19.41'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.30%  '
$ws.Range('E40').Value = '  +0.82%  '
$ws.Range('D41').Value = '''0.373'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.71%  '
$ws.Range('D42').Value = '''5.43'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.81%  '
$ws.Range('D43').Value = '0.0₆0349'
$ws.Range('E43').Value = '  +15.80%  '
$ws.Range('E44').Value = '  -1.16%  '
$ws.Range('E45').Value = '  +3.50%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = '''40.53'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('D48').Value = '''157.54'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('E49').Value = '  +2.72%  '
$ws.Range('E50').Value = '  +1.86%  '
$ws.Range('D51').Value = '''21.96'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.24%  '
